$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.960542333333334
$ws.Range("H2").Value = 14.881627
$ws.Range("I2").Value = 0.1435881646191863
$ws.Range("J2").Value = 0.1435881646191863
$ws.Range("O2").Value = 0.01308107600943097
$ws.Range("P2").Value = 0.01308107600943097
$ws.Range("Q2").Value = 0.4163399715507778
$ws.Range("R2").Value = 3.747059743957
$ws.Range("S2").Value = 0.001878287695438263
$ws.Range("T2").Value = 0.001878287695438263
$ws.Range("G3").Value = 4.960542333333334
$ws.Range("H3").Value = 14.881627
$ws.Range("I3").Value = 0.1435881646191863
$ws.Range("J3").Value = 0.1435881646191863
$ws.Range("M3").Value = 2.908847666666666
$ws.Range("N3").Value = 8.726542999999999
$ws.Range("O3").Value = 0.4533624008902931
$ws.Range("P3").Value = 0.4533624008902929
$ws.Range("Q3").Value = 14.42946199171789
$ws.Range("R3").Value = 129.865157925461
$ws.Range("S3").Value = 0.06509747505118493
$ws.Range("T3").Value = 0.06509747505118492
$ws.Range("G4").Value = 4.960542333333334
$ws.Range("H4").Value = 14.881627
$ws.Range("I4").Value = 0.1435881646191863
$ws.Range("J4").Value = 0.1435881646191863
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.6090526666666666
$ws.Range("N4").Value = 1.827158
$ws.Range("O4").Value = 0.09492472995158634
$ws.Range("P4").Value = 0.09492472995158631
$ws.Range("Q4").Value = 3.021231536229555
$ws.Range("R4").Value = 27.191083826066
$ws.Range("S4").Value = 0.01363006775072018
$ws.Range("T4").Value = 0.01363006775072018
$ws.Range("G5").Value = 4.960542333333334
$ws.Range("H5").Value = 14.881627
$ws.Range("I5").Value = 0.1435881646191863
$ws.Range("J5").Value = 0.1435881646191863
$ws.Range("M5").Value = 2.814333666666667
$ws.Range("N5").Value = 8.443001000000001
$ws.Range("O5").Value = 0.4386317931486897
$ws.Range("P5").Value = 0.4386317931486896
$ws.Range("Q5").Value = 13.96062129362522
$ws.Range("R5").Value = 125.645591642627
$ws.Range("S5").Value = 0.06298233412184291
$ws.Range("T5").Value = 0.06298233412184291
$ws.Range("I6").Value = 0.1393077946862016
$ws.Range("J6").Value = 0.1393077946862016
$ws.Range("O6").Value = 0.01308107600943097
$ws.Range("P6").Value = 0.01308107600943097
$ws.Range("S6").Value = 0.001822295850996408
$ws.Range("T6").Value = 0.001822295850996407
$ws.Range("I7").Value = 0.1393077946862016
$ws.Range("J7").Value = 0.1393077946862016
$ws.Range("M7").Value = 2.908847666666666
$ws.Range("N7").Value = 8.726542999999999
$ws.Range("O7").Value = 0.4533624008902931
$ws.Range("P7").Value = 0.4533624008902929
$ws.Range("Q7").Value = 13.99931905185722
$ws.Range("R7").Value = 125.993871466715
$ws.Range("S7").Value = 0.06315691626166839
$ws.Range("T7").Value = 0.06315691626166836
$ws.Range("I8").Value = 0.1393077946862016
$ws.Range("J8").Value = 0.1393077946862016
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.6090526666666666
$ws.Range("N8").Value = 1.827158
$ws.Range("O8").Value = 0.09492472995158634
$ws.Range("P8").Value = 0.09492472995158631
$ws.Range("Q8").Value = 2.931168482198889
$ws.Range("R8").Value = 26.38051633979
$ws.Range("S8").Value = 0.01322375479073872
$ws.Range("T8").Value = 0.01322375479073872
$ws.Range("I9").Value = 0.1393077946862016
$ws.Range("J9").Value = 0.1393077946862016
$ws.Range("M9").Value = 2.814333666666667
$ws.Range("N9").Value = 8.443001000000001
$ws.Range("O9").Value = 0.4386317931486897
$ws.Range("P9").Value = 0.4386317931486896
$ws.Range("Q9").Value = 13.54445451700056
$ws.Range("R9").Value = 121.900090653005
$ws.Range("S9").Value = 0.06110482778279813
$ws.Range("T9").Value = 0.06110482778279811
$ws.Range("G10").Value = 11.78248366666667
$ws.Range("H10").Value = 35.347451
$ws.Range("I10").Value = 0.3410564996056291
$ws.Range("J10").Value = 0.3410564996056292
$ws.Range("O10").Value = 0.01308107600943097
$ws.Range("P10").Value = 0.01308107600943097
$ws.Range("Q10").Value = 0.9889077816378887
$ws.Range("R10").Value = 8.900170034740999
$ws.Range("S10").Value = 0.004461385994851698
$ws.Range("T10").Value = 0.004461385994851699
$ws.Range("G11").Value = 11.78248366666667
$ws.Range("H11").Value = 35.347451
$ws.Range("I11").Value = 0.3410564996056291
$ws.Range("J11").Value = 0.3410564996056292
$ws.Range("M11").Value = 2.908847666666666
$ws.Range("N11").Value = 8.726542999999999
$ws.Range("O11").Value = 0.4533624008902931
$ws.Range("P11").Value = 0.4533624008902929
$ws.Range("Q11").Value = 34.27345012132144
$ws.Range("R11").Value = 308.461051091893
$ws.Range("S11").Value = 0.1546221935004473
$ws.Range("T11").Value = 0.1546221935004473
$ws.Range("G12").Value = 11.78248366666667
$ws.Range("H12").Value = 35.347451
$ws.Range("I12").Value = 0.3410564996056291
$ws.Range("J12").Value = 0.3410564996056292
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.6090526666666666
$ws.Range("N12").Value = 1.827158
$ws.Range("O12").Value = 0.09492472995158634
$ws.Range("P12").Value = 0.09492472995158631
$ws.Range("Q12").Value = 7.176153097139776
$ws.Range("R12").Value = 64.585377874258
$ws.Range("S12").Value = 0.03237469612329766
$ws.Range("T12").Value = 0.03237469612329766
$ws.Range("G13").Value = 11.78248366666667
$ws.Range("H13").Value = 35.347451
$ws.Range("I13").Value = 0.3410564996056291
$ws.Range("J13").Value = 0.3410564996056292
$ws.Range("M13").Value = 2.814333666666667
$ws.Range("N13").Value = 8.443001000000001
$ws.Range("O13").Value = 0.4386317931486897
$ws.Range("P13").Value = 0.4386317931486896
$ws.Range("Q13").Value = 33.15984046005011
$ws.Range("R13").Value = 298.438564140451
$ws.Range("S13").Value = 0.1495982239870325
$ws.Range("T13").Value = 0.1495982239870325
$ws.Range("G14").Value = 0.9139316666666666
$ws.Range("H14").Value = 2.741795
$ws.Range("I14").Value = 0.02645472244480135
$ws.Range("J14").Value = 0.02645472244480136
$ws.Range("O14").Value = 0.01308107600943097
$ws.Range("P14").Value = 0.01308107600943097
$ws.Range("Q14").Value = 0.07670658942722221
$ws.Range("R14").Value = 0.6903593048449999
$ws.Range("S14").Value = 0.000346056235108846
$ws.Range("T14").Value = 0.000346056235108846
$ws.Range("G15").Value = 0.9139316666666666
$ws.Range("H15").Value = 2.741795
$ws.Range("I15").Value = 0.02645472244480135
$ws.Range("J15").Value = 0.02645472244480136
$ws.Range("M15").Value = 2.908847666666666
$ws.Range("N15").Value = 8.726542999999999
$ws.Range("O15").Value = 0.4533624008902931
$ws.Range("P15").Value = 0.4533624008902929
$ws.Range("Q15").Value = 2.658487996076111
$ws.Range("R15").Value = 23.926391964685
$ws.Range("S15").Value = 0.01199357648246146
$ws.Range("T15").Value = 0.01199357648246146
$ws.Range("G16").Value = 0.9139316666666666
$ws.Range("H16").Value = 2.741795
$ws.Range("I16").Value = 0.02645472244480135
$ws.Range("J16").Value = 0.02645472244480136
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.6090526666666666
$ws.Range("N16").Value = 1.827158
$ws.Range("O16").Value = 0.09492472995158634
$ws.Range("P16").Value = 0.09492472995158631
$ws.Range("Q16").Value = 0.5566325187344443
$ws.Range("R16").Value = 5.009692668609999
$ws.Range("S16").Value = 0.002511207384016938
$ws.Range("T16").Value = 0.002511207384016938
$ws.Range("G17").Value = 0.9139316666666666
$ws.Range("H17").Value = 2.741795
$ws.Range("I17").Value = 0.02645472244480135
$ws.Range("J17").Value = 0.02645472244480136
$ws.Range("M17").Value = 2.814333666666667
$ws.Range("N17").Value = 8.443001000000001
$ws.Range("O17").Value = 0.4386317931486897
$ws.Range("P17").Value = 0.4386317931486896
$ws.Range("Q17").Value = 2.572108658532778
$ws.Range("R17").Value = 23.148977926795
$ws.Range("S17").Value = 0.0116038823432141
$ws.Range("T17").Value = 0.0116038823432141
$ws.Range("G18").Value = 12.077388
$ws.Range("H18").Value = 36.232164
$ws.Range("I18").Value = 0.3495928186441815
$ws.Range("J18").Value = 0.3495928186441815
$ws.Range("O18").Value = 0.01308107600943097
$ws.Range("P18").Value = 0.01308107600943097
$ws.Range("Q18").Value = 1.013659200636
$ws.Range("R18").Value = 9.122932805723998
$ws.Range("S18").Value = 0.004573050233035755
$ws.Range("T18").Value = 0.004573050233035754
$ws.Range("G19").Value = 12.077388
$ws.Range("H19").Value = 36.232164
$ws.Range("I19").Value = 0.3495928186441815
$ws.Range("J19").Value = 0.3495928186441815
$ws.Range("M19").Value = 2.908847666666666
$ws.Range("N19").Value = 8.726542999999999
$ws.Range("O19").Value = 0.4533624008902931
$ws.Range("P19").Value = 0.4533624008902929
$ws.Range("Q19").Value = 35.13128190322799
$ws.Range("R19").Value = 316.1815371290519
$ws.Range("S19").Value = 0.1584922395945309
$ws.Range("T19").Value = 0.1584922395945309
$ws.Range("G20").Value = 12.077388
$ws.Range("H20").Value = 36.232164
$ws.Range("I20").Value = 0.3495928186441815
$ws.Range("J20").Value = 0.3495928186441815
$ws.Range("K20").Value = 3
$ws.Range("L20").Value = 1
$ws.Range("M20").Value = 0.6090526666666666
$ws.Range("N20").Value = 1.827158
$ws.Range("O20").Value = 0.09492472995158634
$ws.Range("P20").Value = 0.09492472995158631
$ws.Range("Q20").Value = 7.355765367767998
$ws.Range("R20").Value = 66.20188830991199
$ws.Range("S20").Value = 0.03318500390281282
$ws.Range("T20").Value = 0.03318500390281282
$ws.Range("G21").Value = 12.077388
$ws.Range("H21").Value = 36.232164
$ws.Range("I21").Value = 0.3495928186441815
$ws.Range("J21").Value = 0.3495928186441815
$ws.Range("M21").Value = 2.814333666666667
$ws.Range("N21").Value = 8.443001000000001
$ws.Range("O21").Value = 0.4386317931486897
$ws.Range("P21").Value = 0.4386317931486896
$ws.Range("Q21").Value = 33.989799653796
$ws.Range("R21").Value = 305.908196884164
$ws.Range("S21").Value = 0.153342524913802
$ws.Range("T21").Value = 0.153342524913802
